# Rename the original sheet to "Rubric"
$wb = $excel.ActiveWorkbook
$rubric = $wb.Worksheets.Item(1)
$rubric.Name = "Rubric"

# Apply bold formatting to the header / total rows (shared by both sheets)
$rubric.Range("A1").Font.Bold = $true
$rubric.Range("A3:C3").Font.Bold = $true
$rubric.Range("A9").Font.Bold = $true

# Duplicate the Rubric sheet (post-bold, pre row-insert) to create the "Grading" sheet
$rubric.Copy($null, $rubric)
$grading = $wb.Worksheets.Item(2)
$grading.Name = "Grading"

# ---- Rubric sheet: insert detail rows under each criterion ----
# Insert bottom-up (using the original, pre-shift row numbers) so the row
# numbers referenced below don't need to account for earlier inserts, and so
# the new shared-strings end up registered in the same order Excel produced
# them (9="...10 to 15 pages", 10="...4 to 6 classes...", 11="...rich media...").
$rubric.Range("A7").EntireRow.Insert()
$rubric.Range("A7").Value = "     10 to 15 pages"

$rubric.Range("A6").EntireRow.Insert()
$rubric.Range("A6").Value = "     4 to 6 classes, 12 to 30 fields"

$rubric.Range("A5").EntireRow.Insert()
$rubric.Range("A5").Value = "     rich media, data, authorization"

# Total row moved from row 9 to row 12 because of the 3 inserted rows; re-bold it
$rubric.Range("A12").Font.Bold = $true

$rubric.Range("B29").Select()

# ---- Grading sheet: update actuals, add notes + percentage column ----
$grading.Range("C5").Value = 5
$grading.Range("D5").Value = "No diagram"

$grading.Range("C7").Value = 14
$grading.Range("D7").Value = "minimal detail"

$grading.Range("F4:F6").Formula = "=C4/B4"
$grading.Range("F7").Formula = "=C7/B7"
$grading.Range("F9").Formula = "=C9/B9"

$grading.Range("F4:F9").NumberFormat = "0.00%"

$grading.Range("A3:D9").Select()
$grading.Activate()
